$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 57: "Longest Palindromic Substrings" (Medium / Aton / Python) ---
# Clone formatting from row 53, which already has the "Medium" (s=3) look.
$ws.Range("A53:I53").Copy()
$ws.Range("A57:I57").PasteSpecial(-4122)
$ws.Range("A57").Value = 5
$ws.Range("B57").Value = "Longest Palindromic Substrings"
$ws.Range("C57").Value = "Dynamic Programming"
$ws.Range("D57").Value = "Aton"
$ws.Range("F57").Value = "Medium"
$ws.Range("G57").Value = "Python"
$ws.Range("E57").Clear()
$ws.Range("H57").Clear()
$ws.Range("I57").Clear()

# --- Row 58: "Best Time to Buy and Sell Stock" (Easy / Aton / Python) ---
# Clone formatting from row 51, which already has the "Easy" (s=6) look.
$ws.Range("A51:I51").Copy()
$ws.Range("A58:I58").PasteSpecial(-4122)
$ws.Range("A58").Value = 121
$ws.Range("B58").Value = "Best Time to Buy and Sell Stock"
$ws.Range("C58").Value = "Dynamic Programming"
$ws.Range("D58").Value = "Aton"
$ws.Range("F58").Value = "Easy"
$ws.Range("G58").Value = "Python"
$ws.Range("E58").Clear()
$ws.Range("H58").Clear()
$ws.Range("I58").Clear()

# --- Row 59: "Range Sum Query - Immutable" (Easy / Aton / Python) ---
$ws.Range("A51:I51").Copy()
$ws.Range("A59:I59").PasteSpecial(-4122)
$ws.Range("A59").Value = 303
$ws.Range("B59").Value = "Range Sum Query - Immutable"
$ws.Range("C59").Value = "Dynamic Programming"
$ws.Range("D59").Value = "Aton"
$ws.Range("F59").Value = "Easy"
$ws.Range("G59").Value = "Python"
$ws.Range("E59").Clear()
$ws.Range("H59").Clear()
$ws.Range("I59").Clear()

# --- Reflect the author's final scroll position / selection ---
$ws.Range("I52").Select() | Out-Null
